# Update column G ("K" - strikeouts) values on the active worksheet to
# reflect the regenerated save_data (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 1
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 2
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
